$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) SmartArt diagram on the slide that asks "Le mie valutazioni sono
#    influenzate dalla nostalgia?" -- the answer node's text is extended
#    from "... le mie valutazi" to "... le mie valutazioni sono aumentate".
#    (This node backs both ppt/diagrams/data2.xml and ppt/diagrams/drawing2.xml;
#    PowerPoint keeps the cached drawing in sync automatically.)
# ---------------------------------------------------------------------------
$targetFragment = "No, anzi col passare del tempo le mie valutazi"
$found = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if (-not $sh.HasSmartArt) { continue }

        $sa = $sh.SmartArt
        $nodes = $sa.AllNodes
        for ($ni = 1; $ni -le $nodes.Count; $ni++) {
            $node = $nodes.Item($ni)
            $tr = $node.TextFrame.TextRange
            $cur = $tr.Text
            if ($cur -and $cur.Contains($targetFragment)) {
                $newText = $cur.Replace(
                    "No, anzi col passare del tempo le mie valutazi",
                    "No, anzi col passare del tempo le mie valutazioni sono aumentate"
                )
                $tr.Text = $newText
                $found = $true
            }
        }
    }
}

if (-not $found) {
    throw "could not locate the SmartArt node containing the 'valutazi' text"
}

# ---------------------------------------------------------------------------
# 2) Plain textbox: the word "aspetti" becomes "fattori" inside the
#    "CasellaDiTesto 5" textbox (the run-level formatting/err flag on that
#    single word is preserved because we only touch that substring).
# ---------------------------------------------------------------------------
$wordFound = $false

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $sl = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $sl.Shapes.Count; $shi++) {
        $sh = $sl.Shapes.Item($shi)
        if (-not $sh.HasTextFrame) { continue }

        $tr = $sh.TextFrame.TextRange
        $full = $tr.Text
        if (-not $full) { continue }

        $idx = $full.IndexOf("vari aspetti")
        if ($idx -ge 0) {
            $wordIdx = $full.IndexOf("aspetti", $idx)
            $sub = $tr.Characters($wordIdx + 1, 7)
            $sub.Text = "fattori"
            $wordFound = $true
        }
    }
}

if (-not $wordFound) {
    throw "could not locate the 'aspetti' word to rename to 'fattori'"
}

Write-Host "Edits applied."
